# Replace formula-driven percentage cells on the "PoFDCtAE" sheet with a
# plain literal value of 1 (100%). The dependent "1 - x" cells (R10, R11,
# R14, R19, R20) keep their existing formulas and simply recalculate to 0
# once their precedent becomes 1.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("PoFDCtAE")

$ws.Range("C3").Value = 1
$ws.Range("D4").Value = 1
$ws.Range("I9").Value = 1
$ws.Range("J10").Value = 1
$ws.Range("K11").Value = 1
$ws.Range("L12").Value = 1
$ws.Range("M13").Value = 1
$ws.Range("N14").Value = 1
$ws.Range("S19").Value = 1
$ws.Range("T20").Value = 1
